$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$timestamp = "2026-02-06 07:00:47"

for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Value = $timestamp
}
